$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 76
$ws1.Range("F4").Value = 98
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 2699
$ws1.Range("F9").Value = 278
$ws1.Range("F10").Value = 135
$ws1.Range("F11").Value = 10178
$ws1.Range("F13").Value = 265
$ws1.Range("F16").Value = 11804
$ws1.Range("F17").Value = 12197
$ws1.Range("F18").Value = 27
$ws1.Range("F19").Value = 98

# Sheet "演出" (sheet2.xml)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 8

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 76
$ws4.Range("F4").Value = 98
$ws4.Range("F6").Value = 50
$ws4.Range("F7").Value = 2699
$ws4.Range("F8").Value = 8
$ws4.Range("F10").Value = 278
$ws4.Range("F11").Value = 135
$ws4.Range("F12").Value = 10178
$ws4.Range("F14").Value = 265
$ws4.Range("F17").Value = 11804
$ws4.Range("F18").Value = 12197
$ws4.Range("F19").Value = 27
$ws4.Range("F20").Value = 98
